# Auto-generated edit script for cryptos.xlsx update
# Updates price (D) and volume-1h (E) columns, and renames row 51 from TheSandbox to Frax
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.141.68"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.862.97"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'0.7094"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "'241.38"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.3100"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'0.07628"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "'24.66"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'0.08347"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.869.46"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "'5.195"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "'0.7072"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "'90.99"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "29.208.58"
$ws.Range("D17").Value = "'5.902"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "2.111.51"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'13.07"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'7.866"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'0.9994"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'0.1583"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "'164.06"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'8.954"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'18.38"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'1.324"
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'4.386"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "'4.253"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "'0.05137"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "'0.7993"
$ws.Range("E34").Value = "  +9.79%  "
$ws.Range("D35").Value = "'1.909"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").Value = "'2.685"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'0.01842"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "'2.691"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "1.163.99"
$ws.Range("E40").Value = "  -5.17%  "
$ws.Range("D41").Value = "'6.218"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'0.8897"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'72.82"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'102.25"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "2.009.57"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'0.5195"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "'1.776"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'9.296"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D51").Value = "'1.0000"
$ws.Range("E51").Value = "  -0.45%  "
